{"js": "// Fix the Czech grammar: \"kv\u016fli fonetiky\" -> \"kv\u016fli fonetice\" inside the\n// summary sentence \"... ale to d\u011bl\u00e1 kv\u016fli fonetiky. Sezn\u00e1m\u00ed se zde i s\n// plukovn\u00edkem Pickeringem, kter\u00e9ho pozve k sob\u011b dom\u016f\". The document also\n// contains an unrelated \"profesor fonetiky\" elsewhere, so we anchor the\n// search on \"kv\u016fli fonetiky\" (unique in the document) rather than the\n// bare word, and we only touch that phrase itself, leaving the rest of\n// the sentence (which contains non-breaking spaces) byte-for-byte intact.\n\nconst originalPhrase = \"kv\u016fli fonetiky\";\nconst correctedPhrase = \"kv\u016fli fonetice\";\n\nconst body = context.document.body;\nconst results = body.search(originalPhrase, { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length !== 1) {\n  throw new Error(\n    `Expected exactly one match for \"${originalPhrase}\", found ${results.items.length}.`\n  );\n}\n\n// Replace in place so the run's formatting (rFonts/color/lang) and the\n// untouched remainder of the sentence are preserved exactly.\nresults.items[0].insertText(correctedPhrase, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Fix the Czech grammar: \"kv\u016fli fonetiky\" -> \"kv\u016fli fonetice\" inside the\n# summary sentence \"... ale to d\u011bl\u00e1 kv\u016fli fonetiky. Sezn\u00e1m\u00ed se zde i s\n# plukovn\u00edkem Pickeringem, kter\u00e9ho pozve k sob\u011b dom\u016f\". The document also\n# contains an unrelated \"profesor fonetiky\" elsewhere, so we anchor the\n# search on \"kv\u016fli fonetiky\" (unique in the document) rather than the\n# bare word, and we only touch that phrase itself, leaving the rest of\n# the sentence (which contains non-breaking spaces) byte-for-byte intact.\n\n$d = $word.ActiveDocument\n\n$originalText = \"kv\u016fli fonetiky\"\n$correctedText = \"kv\u016fli fonetice\"\n\n$replaced = $false\n\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    if ($r.Text.Contains($originalText)) {\n        $find = $r.Find\n        $find.ClearFormatting()\n        $find.Replacement.ClearFormatting()\n        $find.Text = $originalText\n        $find.Replacement.Text = $correctedText\n        $find.Forward = $true\n        $find.Wrap = 0\n        $find.Format = $false\n        $find.MatchCase = $true\n        $find.MatchWholeWord = $false\n        $find.MatchWildcards = $false\n        $find.MatchSoundsLike = $false\n        $find.MatchAllWordForms = $false\n\n        $replaced = $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $correctedText, 1)\n        if ($replaced) {\n            break\n        }\n    }\n}\n\nif (-not $replaced) {\n    throw \"Could not find the target sentence to correct.\"\n}\n"}
